$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B. This shifts the existing
# B:K columns (PercActivations ... totalStd) one column to the right,
# becoming C:L, and leaves a fresh empty column B behind.
$ws.Range("B1").EntireColumn.Insert()

# Give the new B1 header the same look (bold / bordered / centered) as
# the other header cells, then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value = "segments"

# The segment-name labels used to live in column A (with the bold/boxed
# header-ish style). Move each one into the new column B as plain text,
# and replace column A's content with a 0-based numeric row index.
$lastRow = 20
for ($r = 2; $r -le $lastRow; $r++) {
    $segmentName = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $segmentName
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column B's data cells should be plain (unstyled) text, not carry over
# the bordered/bold style that Insert copied in from column A.
$ws.Range("B2:B" + $lastRow).ClearFormats()
